$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "warning" cells that used to hold placeholder default-object-type
# names; they now stay empty and get flagged with a red fill instead.
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()

# Highlight the now-empty cells in red to visually warn the user.
$ws.Range("D5:D6").Interior.ColorIndex = 3

# Match the updated selection recorded by Excel after this edit.
$ws.Range("D5:D6").Select() | Out-Null
